# Revert "Mudança do Sprint 002"
# 1) Remove the picture-only slide (4th slide, sldId 282) that the
#    reverted commit had inserted between "1a Estoria" and "2a Estoria".
# 2) Roll the footer date placeholders (master + every layout) back
#    from 9/17/2018 to 9/16/2018.

$p = $ppt.ActivePresentation

# --- 1) delete the inserted picture slide -------------------------------
# It is the 4th slide in the deck (after title, "Estorias do Sprint
# Atual" and "1a Estoria"), and before "2a Estoria".
$p.Slides.Item(4).Delete()

# --- 2) fix the cached date field text ----------------------------------
function Update-DatePlaceholders($shapes, $newText) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.Name -like "Date Placeholder*") {
            if ($sh.HasTextFrame -and $sh.TextFrame.HasText) {
                $sh.TextFrame.TextRange.Text = $newText
            }
        }
    }
}

$master = $p.SlideMaster
Update-DatePlaceholders $master.Shapes "9/16/2018"

for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Update-DatePlaceholders $layout.Shapes "9/16/2018"
}
